$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.834.22"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "3.411.50"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.26"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.18"
$ws.Range("E6").Value = "  +2.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -1.38%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.726"
$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("E10").Value = "  -5.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.87"
$ws.Range("E11").Value = "  +1.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.12"
$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("D13").Value = "3.950.30"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000208"
$ws.Range("E15").Value = "  -1.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.45"
$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").Value = "3.400.79"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("E18").Value = "  +1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.29"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").Value = "61.884.52"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.14"
$ws.Range("E21").Value = "  +18.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.83"
$ws.Range("E22").Value = "  +0.86%  "

$ws.Range("E23").Value = "  +3.08%  "

$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("E25").Value = "  +2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.66"
$ws.Range("E26").Value = "  +13.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.01"
$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.75"
$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  +6.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.75"
$ws.Range("E30").Value = "  +0.96%  "

$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("E32").Value = "  -1.83%  "

$ws.Range("E33").Value = "  -4.82%  "

$ws.Range("E34").Value = "  -4.13%  "

$ws.Range("E35").Value = "  -0.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.78"
$ws.Range("E36").Value = "  +5.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0487"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.06"
$ws.Range("E39").Value = "  +5.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.332"
$ws.Range("E40").Value = "  +7.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.11"
$ws.Range("E41").Value = "  +5.96%  "

$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("E44").Value = "  +5.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.60"
$ws.Range("E45").Value = "  +8.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.25"
$ws.Range("E46").Value = "  +5.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.56"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  +19.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.91"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("B50").Value = "Fetch.AI"
$ws.Range("C50").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.96"
$ws.Range("E50").Value = "  +19.76%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.142"
$ws.Range("E51").Value = "  +8.59%  "
